# Weekly Fruta/Hortalizas update: a new week's price record is added for
# "Feria Lagunitas de Puerto Montt - Mango", inserted at row 174. This
# pushes the existing rows 174-203 down to 175-204 (the last existing row,
# old 203, becomes new row 204) and grows the used range from A1:T203 to
# A1:T204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 174, shifting rows 174:203
# down to 175:204 (formats are carried along automatically, matching the
# existing "s=2" date style already used by column D).
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly record.
$ws.Cells.Item(174, 1).Value  = 4
$ws.Cells.Item(174, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(174, 3).Value  = "Los Lagos"
$ws.Cells.Item(174, 4).Value  = 44694
$ws.Cells.Item(174, 5).Value  = 10
$ws.Cells.Item(174, 6).Value  = "Fruta"
$ws.Cells.Item(174, 7).Value  = 100108
$ws.Cells.Item(174, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(174, 9).Value  = 100108002
$ws.Cells.Item(174, 10).Value = "Mango"
$ws.Cells.Item(174, 11).Value = "Sin especificar"
$ws.Cells.Item(174, 12).Value = "Primera"
$ws.Cells.Item(174, 13).Value = 160
$ws.Cells.Item(174, 14).Value = 7500
$ws.Cells.Item(174, 15).Value = 8000
$ws.Cells.Item(174, 16).Value = 7750
$ws.Cells.Item(174, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(174, 18).Value = "Perú"
$ws.Cells.Item(174, 19).Value = 1938
$ws.Cells.Item(174, 20).Value = 4
